$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column D. This shifts the existing D:J
#    (fiscal-year) data one column to the right (D->E, E->F, ... J->K),
#    matching the workbook gaining a new (most recent) fiscal year of data.
$ws.Columns("D:D").Insert()

# 2. The newly inserted column D cells have no number formatting yet.
#    Copy the formatting from column E (which holds what used to be column D)
#    into column D for every row that actually carries data, so the new
#    column visually matches (date format for the "Period Ending" rows,
#    numeric format for the data rows).
$fmtSrc = "E7:E10,E12:E15,E17:E18,E20:E35,E38:E38,E41:E54,E57:E66,E68:E77,E80:E81,E83:E89,E91:E94,E96:E102"
$fmtDst = "D7:D10,D12:D15,D17:D18,D20:D35,D38:D38,D41:D54,D57:D66,D68:D77,D80:D81,D83:D89,D91:D94,D96:D102"
$ws.Range($fmtSrc).Copy()
$ws.Range($fmtDst).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column D should keep the same width/bestFit behaviour as the other
# (now shifted) data columns E:J.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# 3. Populate the new column D with the new fiscal year's figures (and the
#    new period-ending date 2018-12-31 / serial 43465 on the three
#    "Period Ending" header rows). Rows that previously had no disclosed
#    breakout (value 0 across the whole row) simply continue to read 0 in
#    the new column.
$newColumnD = @{
    7 = 43465;
    8 = 1073200;
    9 = 477000;
    10 = 596300;
    12 = 225700;
    13 = 0;
    14 = 0;
    15 = 11400;
    17 = 1165100;
    18 = -91900;
    20 = 27400;
    21 = -37500;
    22 = 0;
    23 = -64600;
    24 = 0;
    25 = 0;
    26 = -64600;
    27 = -64600;
    28 = 0;
    29 = 0;
    30 = 0;
    31 = 0;
    32 = -27400;
    33 = -64600;
    34 = 0;
    35 = -64600;
    38 = 43465;
    41 = 502600;
    42 = 1559000;
    43 = 41300;
    44 = 0;
    45 = 26200;
    46 = 2129100;
    47 = 0;
    48 = 61600;
    49 = 64100;
    50 = 0;
    51 = 0;
    52 = 0;
    53 = 0;
    54 = 2254800;
    57 = 97000;
    58 = 2600;
    59 = 39200;
    60 = 138700;
    61 = 22300;
    62 = 3000;
    63 = 0;
    64 = 0;
    65 = 0;
    66 = 164000;
    68 = 0;
    69 = 0;
    70 = 0;
    71 = 0;
    72 = -187800;
    73 = 0;
    74 = 0;
    75 = 0;
    76 = 2090800;
    77 = 0;
    80 = 43465;
    81 = -64600;
    83 = 27100;
    84 = 0;
    85 = 0;
    86 = 0;
    87 = 0;
    88 = 0;
    89 = 9300;
    91 = -28000;
    92 = 0;
    93 = 0;
    94 = -810600;
    96 = 0;
    97 = 0;
    98 = 0;
    99 = 0;
    100 = 1072200;
    101 = -1900;
    102 = 269000;

}
foreach ($r in $newColumnD.Keys) {
    $addr = "D" + $r
    $ws.Range($addr).Value = $newColumnD[$r]
}

# 4. Row 91 ("Capital Expenditures") also received corrected historical
#    figures for the years that used to sit in D:I (now E:J), not just a
#    straight shift -- update those explicitly.
$ws.Range("E91").Value = -20000
$ws.Range("F91").Value = -23800
$ws.Range("G91").Value = -16500
$ws.Range("H91").Value = -20600
$ws.Range("I91").Value = -3500
$ws.Range("J91").Value = -1700
